# Update the build-version timestamp embedded throughout the workbook.
$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldText)) {
            $cell.Value = $val.Replace($oldText, $newText)
        }
    }
}
